$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 12.31540361230704
$ws.Range("D2").Value = 4.231818227149851
$ws.Range("E2").Value = 12.62927869702264
$ws.Range("F2").Value = 25.12074969784221
$ws.Range("G2").Value = 31.09846476813864
$ws.Range("H2").Value = 14.19147225454205
$ws.Range("L2").Value = 8.89678928261948
$ws.Range("M2").Value = 22.58615359304756
$ws.Range("N2").Value = 17.38958903131407
$ws.Range("O2").Value = 22.06309919922899

$ws.Range("C3").Value = 12.35085028893442
$ws.Range("D3").Value = 4.231742074520271
$ws.Range("E3").Value = 12.71486225660265
$ws.Range("F3").Value = 24.89328269715379
$ws.Range("G3").Value = 30.60324944919293
$ws.Range("H3").Value = 14.18494371350316
$ws.Range("L3").Value = 8.935730736436721
$ws.Range("M3").Value = 21.88444132594487
$ws.Range("N3").Value = 17.10576821282844
$ws.Range("O3").Value = 21.95935382316665

$ws.Range("C4").Value = 12.37616480145269
$ws.Range("D4").Value = 4.231754909150133
$ws.Range("E4").Value = 12.77018974641919
$ws.Range("F4").Value = 24.76085180009394
$ws.Range("G4").Value = 30.30703916771309
$ws.Range("H4").Value = 14.18428144291796
$ws.Range("L4").Value = 8.960784555437689
$ws.Range("M4").Value = 21.44211288843011
$ws.Range("N4").Value = 16.9310585898549
$ws.Range("O4").Value = 21.90201447210531

$ws.Range("C5").Value = 12.38736878506425
$ws.Range("D5").Value = 4.231775323363022
$ws.Range("E5").Value = 12.79343537064705
$ws.Range("F5").Value = 24.70876159359317
$ws.Range("G5").Value = 30.18848069176239
$ws.Range("H5").Value = 14.18485323700378
$ws.Range("L5").Value = 8.971282657620701
$ws.Range("M5").Value = 21.25921551990863
$ws.Range("N5").Value = 16.85983833098167
$ws.Range("O5").Value = 21.88026505771417

$ws.Range("C6").Value = 12.38928269944944
$ws.Range("D6").Value = 4.231779636046952
$ws.Range("E6").Value = 12.79733753459915
$ws.Range("F6").Value = 24.70022697709681
$ws.Range("G6").Value = 30.16892899146946
$ws.Range("H6").Value = 14.18499901383839
$ws.Range("L6").Value = 8.973043308587108
$ws.Range("M6").Value = 21.22869322257825
$ws.Range("N6").Value = 16.84801332934872
$ws.Range("O6").Value = 21.87675169439862

$ws.Range("C7").Value = 12.37631231258949
$ws.Range("D7").Value = 4.231755122733434
$ws.Range("E7").Value = 12.77050041328918
$ws.Range("F7").Value = 24.76014162324568
$ws.Range("G7").Value = 30.30543131562927
$ws.Range("H7").Value = 14.18428574658131
$ws.Range("L7").Value = 8.960924967346109
$ws.Range("M7").Value = 21.4396566474777
$ws.Range("N7").Value = 16.93009807297626
$ws.Range("O7").Value = 21.90171458440284

$ws.Range("C8").Value = 12.32688610704167
$ws.Range("D8").Value = 4.231779719948024
$ws.Range("E8").Value = 12.65821143292973
$ws.Range("F8").Value = 25.04085309656326
$ws.Range("G8").Value = 30.92619054708625
$ws.Range("H8").Value = 14.18852698777637
$ws.Range("L8").Value = 8.909979482816221
$ws.Range("M8").Value = 22.34672265466789
$ws.Range("N8").Value = 17.29187093185413
$ws.Range("O8").Value = 22.02601882417525

$ws.Range("C9").Value = 12.25832547243187
$ws.Range("D9").Value = 4.232292195346705
$ws.Range("E9").Value = 12.46004043767267
$ws.Range("F9").Value = 25.64602259496743
$ws.Range("G9").Value = 32.19796156294807
$ws.Range("H9").Value = 14.223363201697
$ws.Range("L9").Value = 8.819107992887147
$ws.Range("M9").Value = 24.02489547832061
$ws.Range("N9").Value = 17.99407867198627
$ws.Range("O9").Value = 22.31936591299453

$ws.Range("C10").Value = 12.2254977589977
$ws.Range("D10").Value = 4.23294036765294
$ws.Range("E10").Value = 12.32784228399085
$ws.Range("F10").Value = 26.12011582875498
$ws.Range("G10").Value = 33.154916789998
$ws.Range("H10").Value = 14.26504844224996
$ws.Range("L10").Value = 8.757792076337637
$ws.Range("M10").Value = 25.18516056509593
$ws.Range("N10").Value = 18.50076743103616
$ws.Range("O10").Value = 22.56377152927111

$ws.Range("C11").Value = 12.21442151483202
$ws.Range("D11").Value = 4.233291842728237
$ws.Range("E11").Value = 12.27060577515923
$ws.Range("F11").Value = 26.34131567546637
$ws.Range("G11").Value = 33.59293277453217
$ws.Range("H11").Value = 14.28747812539641
$ws.Range("L11").Value = 8.731068080148068
$ws.Range("M11").Value = 25.69518991400071
$ws.Range("N11").Value = 18.72830678788785
$ws.Range("O11").Value = 22.68089635239924

$ws.Range("C12").Value = 12.21078538766895
$ws.Range("D12").Value = 4.233432879809946
$ws.Range("E12").Value = 12.2493487983609
$ws.Range("F12").Value = 26.42579775361004
$ws.Range("G12").Value = 33.75900383266122
$ws.Range("H12").Value = 14.29646681077201
$ws.Range("L12").Value = 8.721115542550878
$ws.Range("M12").Value = 25.88562315376026
$ws.Range("N12").Value = 18.81396908886077
$ws.Range("O12").Value = 22.72607204383014

$ws.Range("C13").Value = 12.2115436134451
$ws.Range("D13").Value = 4.233402154944031
$ws.Range("E13").Value = 12.25390830956689
$ws.Range("F13").Value = 26.40757230806879
$ws.Range("G13").Value = 33.72323124847918
$ws.Range("H13").Value = 14.29450898062282
$ws.Range("L13").Value = 8.723251573829929
$ws.Range("M13").Value = 25.84473236501429
$ws.Range("N13").Value = 18.79554368596318
$ws.Range("O13").Value = 22.71630656611913

$ws.Range("C14").Value = 12.21411116205991
$ws.Range("D14").Value = 4.233303287927324
$ws.Range("E14").Value = 12.26884858796889
$ws.Range("F14").Value = 26.34825211191873
$ws.Range("G14").Value = 33.60659248373955
$ws.Range("H14").Value = 14.28820773039975
$ws.Range("L14").Value = 8.730245931491394
$ws.Range("M14").Value = 25.71091194240824
$ws.Range("N14").Value = 18.73536476286892
$ws.Range("O14").Value = 22.68459667554239

$ws.Range("C15").Value = 12.21575665415226
$ws.Range("D15").Value = 4.233243757258784
$ws.Range("E15").Value = 12.27805428110413
$ws.Range("F15").Value = 26.31200804467581
$ws.Range("G15").Value = 33.53516908261033
$ws.Range("H15").Value = 14.28441238971251
$ws.Range("L15").Value = 8.734551935132066
$ws.Range("M15").Value = 25.62858687644218
$ws.Range("N15").Value = 18.69843585400502
$ws.Range("O15").Value = 22.66527967114886

$ws.Range("C16").Value = 12.22629954509448
$ws.Range("D16").Value = 4.232918520060943
$ws.Range("E16").Value = 12.33164120654313
$ws.Range("F16").Value = 26.10576466218518
$ws.Range("G16").Value = 33.12632980656782
$ws.Range("H16").Value = 14.26365208778559
$ws.Range("L16").Value = 8.759562002988874
$ws.Range("M16").Value = 25.15145820804423
$ws.Range("N16").Value = 18.4858310628749
$ws.Range("O16").Value = 22.55623408395402

$ws.Range("C17").Value = 12.23375781978279
$ws.Range("D17").Value = 4.232733355789524
$ws.Range("E17").Value = 12.36525811956231
$ws.Range("F17").Value = 25.98060354758177
$ws.Range("G17").Value = 32.87607770330126
$ws.Range("H17").Value = 14.25180197219421
$ws.Range("L17").Value = 8.775203666527608
$ws.Range("M17").Value = 24.85409120640824
$ws.Range("N17").Value = 18.35459196480615
$ws.Range("O17").Value = 22.49083850568172

$ws.Range("C18").Value = 12.2384104926922
$ws.Range("D18").Value = 4.232632195235852
$ws.Range("E18").Value = 12.38486673479274
$ws.Range("F18").Value = 25.90913962466716
$ws.Range("G18").Value = 32.73240339167113
$ws.Range("H18").Value = 14.2453126618174
$ws.Range("L18").Value = 8.784310411375452
$ws.Range("M18").Value = 24.68138761568975
$ws.Range("N18").Value = 18.27883284988673
$ws.Range("O18").Value = 22.45378613890925

$ws.Range("C19").Value = 12.24004801811843
$ws.Range("D19").Value = 4.232598867484856
$ws.Range("E19").Value = 12.39155277155633
$ws.Range("F19").Value = 25.8850357434003
$ws.Range("G19").Value = 32.68380876272071
$ws.Range("H19").Value = 14.24317167669901
$ws.Range("L19").Value = 8.787412729349473
$ws.Range("M19").Value = 24.62263176847287
$ws.Range("N19").Value = 18.25313748323822
$ws.Range("O19").Value = 22.44133820559873

$ws.Range("C20").Value = 12.23292629131125
$ws.Range("D20").Value = 4.232752515508862
$ws.Range("E20").Value = 12.36165128117915
$ws.Range("F20").Value = 25.99387333574991
$ws.Range("G20").Value = 32.90269147988026
$ws.Range("H20").Value = 14.25302966251011
$ws.Range("L20").Value = 8.773527198538138
$ws.Range("M20").Value = 24.88591996848558
$ws.Range("N20").Value = 18.36859150265975
$ws.Range("O20").Value = 22.49774209289638

$ws.Range("C21").Value = 12.21334183464102
$ws.Range("D21").Value = 4.233332113638435
$ws.Range("E21").Value = 12.26444894115219
$ws.Range("F21").Value = 26.36565701052886
$ws.Range("G21").Value = 33.64084801522299
$ws.Range("H21").Value = 14.29004515541369
$ws.Range("L21").Value = 8.728186986305825
$ws.Range("M21").Value = 25.75029271392254
$ws.Range("N21").Value = 18.7530549869502
$ws.Range("O21").Value = 22.69388856188441

$ws.Range("C22").Value = 12.20379732295377
$ws.Range("D22").Value = 4.233757109195006
$ws.Range("E22").Value = 12.20335371134536
$ws.Range("F22").Value = 26.61279078814983
$ws.Range("G22").Value = 34.12438365921447
$ws.Range("H22").Value = 14.31712043176291
$ws.Range("L22").Value = 8.699529072972945
$ws.Range("M22").Value = 26.29939367016662
$ws.Range("N22").Value = 19.00136560686324
$ws.Range("O22").Value = 22.82686321122026

$ws.Range("C23").Value = 12.20859248511156
$ws.Range("D23").Value = 4.233526120193635
$ws.Range("E23").Value = 12.23573880389736
$ws.Range("F23").Value = 26.48053664100011
$ws.Range("G23").Value = 33.86626929378798
$ws.Range("H23").Value = 14.30240730587883
$ws.Range("L23").Value = 8.714735443426713
$ws.Range("M23").Value = 26.007819702187
$ws.Range("N23").Value = 18.8691326645061
$ws.Range("O23").Value = 22.75546550393854

$ws.Range("C24").Value = 12.23330108924546
$ws.Range("D24").Value = 4.232743836872693
$ws.Range("E24").Value = 12.3632810549724
$ws.Range("F24").Value = 25.98787252331808
$ws.Range("G24").Value = 32.89065875762552
$ws.Range("H24").Value = 14.25247361577621
$ws.Range("L24").Value = 8.774284774072372
$ws.Range("M24").Value = 24.87153560389041
$ws.Range("N24").Value = 18.36226326508119
$ws.Range("O24").Value = 22.49461928280682

$ws.Range("C25").Value = 12.27380927573925
$ws.Range("D25").Value = 4.232105035291562
$ws.Range("E25").Value = 12.51129613266566
$ws.Range("F25").Value = 25.47684072858808
$ws.Range("G25").Value = 31.84918068100226
$ws.Range("H25").Value = 14.21110553433758
$ws.Range("L25").Value = 8.842730275542047
$ws.Range("M25").Value = 23.58292278396566
$ws.Range("N25").Value = 17.99407867198627
$ws.Range("O25").Value = 22.2348289538332
